# Applies the textual edits described by the diff to the Korean product
# specification document. Each replacement targets a unique run of text,
# so whole-phrase Find/Replace is safe and precise.

$d = $word.ActiveDocument

$replacements = @(
    @{
        Old = "복원력이 우수한 고급 보안 제품인 Contoso CipherGuard Sentinel X7은 컴퓨터 네트워크 인프라를 강화하여 다양한 위협과 취약성을 방지할 수 있도록 세심하게 제작되었습니다. "
        New = "Contoso CipherGuard Sentinel X7은 다양한 위협과 취약성에 대비하여 컴퓨터 네트워크 인프라를 강화하도록 세심하게 설계된 복원력 있는 고급 보안 제품입니다. "
    },
    @{
        Old = " 다계층 방어 접근 방식을 사용하는 엔드포인트 보안 모듈은 바이러스 백신, 맬웨어 방지 및 호스트 기반 침입 방지 기능을 통합합니다. "
        New = " 다중 계층 방어 접근 방식을 사용하는 엔드포인트 보안 모듈은 바이러스 백신, 맬웨어 방지 및 호스트 기반 침입 방지 기능을 통합합니다. "
    },
    @{
        Old = " Windows Server 2019 이상, CentOS 8 또는 해당 버전과 호환"
        New = " Windows Server 2019 이상, CentOS 8 또는 동급 버전과 호환"
    },
    @{
        Old = " Contoso는 Contoso CipherGuard Sentinel X7과 관련된 기술 문제 또는 문의에 대한 신속한 지원을 보장하기 위해 전용 24/7 지원 팀을 제공합니다."
        New = " Contoso는 Contoso CipherGuard Sentinel X7과 관련된 기술 문제 또는 문의에 대한 신속한 지원을 보장하기 위해 연중무휴 상시 전담 지원팀을 제공합니다."
    }
)

foreach ($r in $replacements) {
    $find = $d.Content.Find
    $find.ClearFormatting()
    # Execute(FindText, MatchCase, MatchWholeWord, MatchWildcards, MatchSoundsLike,
    #         MatchAllWordForms, Forward, Wrap, Format, ReplaceWith, Replace)
    $find.Execute($r.Old, $true, $false, $false, $false, $false, $true, 1, $false, $r.New, 2)
}
